$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 currently holds the text "R40" (shared string). It needs to become
# the text "1" - note this must remain a *text* value (shared string), not a
# number, and the cell's existing style (border/fill formatting) must be kept.
$cell = $ws.Range("B11")

# Stash the current formatting on a scratch cell far outside the used range
# so we can restore it after the value/number-format change below (setting
# NumberFormat to force text storage otherwise creates a brand new style).
$scratch = $ws.Range("Z100")
$cell.Copy($scratch)

# Force the new value to be stored as text rather than being auto-converted
# to a number.
$cell.NumberFormat = "@"
$cell.Value = "1"

# Restore the original cell formatting/style.
$scratch.Copy()
$cell.PasteSpecial(-4122)  # xlPasteFormats

# Clean up the scratch cell.
$scratch.Clear()
